$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, pushing existing rows 61:85 down to 62:86
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with the new weekly record
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44636
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112021
$ws.Range("G61").Value = "Ají"
$ws.Range("H61").Value = "Americana (o)"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 25
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 11000
$ws.Range("M61").Value = 10400
$ws.Range("N61").Value = "$/caja 15 kilos"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 693
$ws.Range("Q61").Value = 15
$ws.Range("R61").Value = "Hortaliza"
